$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 3499.923
$ws.Range("J40").Value2 = 3045.182
$ws.Range("L40").Value2 = 3045.182
$ws.Range("N40").Value2 = -3395.182
$ws.Range("H46").Value2 = 19333
$ws.Range("J46").Value2 = 26999.5
$ws.Range("L46").Value2 = 80998.5
$ws.Range("N46").Value2 = -81236.5
$ws.Range("H60").Value2 = 19333
$ws.Range("J60").Value2 = 26999.5
$ws.Range("L60").Value2 = 80998.5
$ws.Range("N60").Value2 = -81966.5
$ws.Range("H70").Value2 = 1286021.1
$ws.Range("J70").Value2 = 1869.75
$ws.Range("L70").Value2 = 5609.25
$ws.Range("N70").Value2 = -6149.25
$ws.Range("H73").Value2 = 1286021.1
$ws.Range("J73").Value2 = 1869.75
$ws.Range("L73").Value2 = 5609.25
$ws.Range("N73").Value2 = -7481.25
$ws.Range("H80").Value2 = 813027.4
$ws.Range("I80").Value2 = 1263999.5
$ws.Range("J80").Value2 = 1277.5
$ws.Range("K80").Value2 = 3791998.5
$ws.Range("L80").Value2 = 3832.5
$ws.Range("M80").Value2 = -3791000.5
$ws.Range("N80").Value2 = -5828.5
$ws.Range("H83").Value2 = 813027.4
$ws.Range("I83").Value2 = 1263999.5
$ws.Range("J83").Value2 = 1277.5
$ws.Range("K83").Value2 = 11375995.5
$ws.Range("L83").Value2 = 11497.5
$ws.Range("M83").Value2 = -11371003.5
$ws.Range("N83").Value2 = -21481.5
$ws.Range("H86").Value2 = 3453075
$ws.Range("I86").Value2 = 2859.6843
$ws.Range("J86").Value2 = 10008484
$ws.Range("K86").Value2 = 2859.6843
$ws.Range("L86").Value2 = 10008484
$ws.Range("M86").Value2 = -1736.6843
$ws.Range("N86").Value2 = -10010730
$ws.Range("H89").Value2 = 3453075
$ws.Range("I89").Value2 = 2859.6843
$ws.Range("J89").Value2 = 10008484
$ws.Range("K89").Value2 = 14298.4215
$ws.Range("L89").Value2 = 50042420
$ws.Range("M89").Value2 = -8682.4215
$ws.Range("N89").Value2 = -50053652
$ws.Range("H132").Value2 = 5675.3335
$ws.Range("I132").Value2 = 5220.25
$ws.Range("J132").Value2 = 7950.75
$ws.Range("K132").Value2 = 15660.75
$ws.Range("L132").Value2 = 23852.25
$ws.Range("M132").Value2 = -13130.75
$ws.Range("N132").Value2 = -28912.25
$ws.Range("H138").Value2 = 2157.681
$ws.Range("J138").Value2 = 2886.5186
$ws.Range("L138").Value2 = 8659.5558
$ws.Range("N138").Value2 = -18939.5558

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 1671.3636
$ws.Range("I32").Value2 = 1671.3636
$ws.Range("K32").Value2 = 1671.3636
$ws.Range("M32").Value2 = -1384.3636
$ws.Range("H61").Value2 = 19610486
$ws.Range("I61").Value2 = 27779564
$ws.Range("K61").Value2 = 27779564
$ws.Range("M61").Value2 = -27779352
$ws.Range("H74").Value2 = 3441.3333
$ws.Range("I74").Value2 = 3798.4
$ws.Range("K74").Value2 = 3798.4
$ws.Range("M74").Value2 = -2924.4
$ws.Range("H77").Value2 = 3441.3333
$ws.Range("I77").Value2 = 3798.4
$ws.Range("K77").Value2 = 18992
$ws.Range("M77").Value2 = -14624
$ws.Range("H102").Value2 = 13305.071
$ws.Range("I102").Value2 = 3296.4546
$ws.Range("K102").Value2 = 3296.4546
$ws.Range("M102").Value2 = -1674.4546
$ws.Range("H136").Value2 = 19610486
$ws.Range("I136").Value2 = 27779564
$ws.Range("K136").Value2 = 83338692
$ws.Range("M136").Value2 = -83336142
$ws.Range("H43").Value2 = 2666.5557
$ws.Range("I43").Value2 = 2666.5557
$ws.Range("K43").Value2 = 2666.5557
$ws.Range("M43").Value2 = -2515.5557

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 216.66667
$ws.Range("I22").Value2 = 232.14285
$ws.Range("J22").Value2 = 162.5
$ws.Range("K22").Value2 = 232.14285
$ws.Range("L22").Value2 = 162.5
$ws.Range("M22").Value2 = -59.14285000000001
$ws.Range("N22").Value2 = -508.5
$ws.Range("H86").Value2 = 4206.4443
$ws.Range("I86").Value2 = 3981.375
$ws.Range("K86").Value2 = 3981.375
$ws.Range("M86").Value2 = -2858.375
$ws.Range("H89").Value2 = 4206.4443
$ws.Range("I89").Value2 = 3981.375
$ws.Range("K89").Value2 = 19906.875
$ws.Range("M89").Value2 = -14290.875
$ws.Range("H99").Value2 = 906.4
$ws.Range("I99").Value2 = 780.5
$ws.Range("J99").Value2 = 1410
$ws.Range("K99").Value2 = 780.5
$ws.Range("L99").Value2 = 1410
$ws.Range("M99").Value2 = 717.5
$ws.Range("N99").Value2 = -4406
$ws.Range("H107").Value2 = 20009476
$ws.Range("I107").Value2 = 11002.286
$ws.Range("J107").Value2 = 125001464
$ws.Range("K107").Value2 = 11002.286
$ws.Range("L107").Value2 = 125001464
$ws.Range("M107").Value2 = -9082.286
$ws.Range("N107").Value2 = -125005304

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value2 = 30000
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 30000
$ws.Range("K42").Value2 = 0
$ws.Range("L42").Value2 = 30000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value2 = -31186
$ws.Range("H62").Value2 = 52637050
$ws.Range("I62").Value2 = 5517.636
$ws.Range("J62").Value2 = 125005400
$ws.Range("K62").Value2 = 5517.636
$ws.Range("L62").Value2 = 125005400
$ws.Range("M62").Value2 = -4893.636
$ws.Range("N62").Value2 = -125006648
$ws.Range("H65").Value2 = 52637050
$ws.Range("I65").Value2 = 5517.636
$ws.Range("J65").Value2 = 125005400
$ws.Range("K65").Value2 = 27588.18
$ws.Range("L65").Value2 = 625027000
$ws.Range("M65").Value2 = -24468.18
$ws.Range("N65").Value2 = -625033240
$ws.Range("H94").Value2 = 1519.25
$ws.Range("J94").Value2 = 1609
$ws.Range("L94").Value2 = 1609
$ws.Range("N94").Value2 = -2511
$ws.Range("H107").Value2 = 1940.875
$ws.Range("I107").Value2 = 2125.6428
$ws.Range("J107").Value2 = 647.5
$ws.Range("K107").Value2 = 2125.6428
$ws.Range("L107").Value2 = 647.5
$ws.Range("M107").Value2 = -205.6428000000001
$ws.Range("N107").Value2 = -4487.5
$ws.Range("H134").Value2 = 3014.7917
$ws.Range("I134").Value2 = 2117.3572
$ws.Range("K134").Value2 = 6352.071599999999
$ws.Range("M134").Value2 = -3817.071599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 811.6429000000001
$ws.Range("J2").Value2 = 192
$ws.Range("L2").Value2 = 1152
$ws.Range("N2").Value2 = -1378
$ws.Range("H76").Value2 = 7576.5
$ws.Range("I76").Value2 = 999
$ws.Range("K76").Value2 = 2997
$ws.Range("M76").Value2 = -2614
$ws.Range("H79").Value2 = 7576.5
$ws.Range("I79").Value2 = 999
$ws.Range("K79").Value2 = 2997
$ws.Range("M79").Value2 = -1671
$ws.Range("H116").Value2 = 2325
$ws.Range("J116").Value2 = 4500
$ws.Range("L116").Value2 = 13500
$ws.Range("N116").Value2 = -20384
$ws.Range("H118").Value2 = 1190
$ws.Range("I118").Value2 = 376
$ws.Range("J118").Value2 = 3632
$ws.Range("K118").Value2 = 1128
$ws.Range("L118").Value2 = 10896
$ws.Range("M118").Value2 = 115
$ws.Range("N118").Value2 = -13382

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 3998
$ws.Range("I80").Value2 = 3997
$ws.Range("K80").Value2 = 3997
$ws.Range("M80").Value2 = -2999
$ws.Range("H83").Value2 = 3998
$ws.Range("I83").Value2 = 3997
$ws.Range("K83").Value2 = 19985
$ws.Range("M83").Value2 = -14993
$ws.Range("H107").Value2 = 631.36365
$ws.Range("I107").Value2 = 443.375
$ws.Range("J107").Value2 = 1132.6666
$ws.Range("K107").Value2 = 443.375
$ws.Range("L107").Value2 = 1132.6666
$ws.Range("M107").Value2 = 1476.625
$ws.Range("N107").Value2 = -4972.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 852.2
$ws.Range("I16").Value2 = 815.375
$ws.Range("J16").Value2 = 999.5
$ws.Range("K16").Value2 = 815.375
$ws.Range("L16").Value2 = 999.5
$ws.Range("M16").Value2 = -645.375
$ws.Range("N16").Value2 = -1339.5
$ws.Range("H46").Value2 = 1703.7
$ws.Range("I46").Value2 = 1000
$ws.Range("J46").Value2 = 1740.7368
$ws.Range("K46").Value2 = 1000
$ws.Range("L46").Value2 = 1740.7368
$ws.Range("M46").Value2 = -812
$ws.Range("N46").Value2 = -2116.7368
$ws.Range("H82").Value2 = 3922.0908
$ws.Range("I82").Value2 = 4250.4443
$ws.Range("J82").Value2 = 2444.5
$ws.Range("K82").Value2 = 4250.4443
$ws.Range("L82").Value2 = 2444.5
$ws.Range("M82").Value2 = -3889.4443
$ws.Range("N82").Value2 = -3166.5
$ws.Range("H85").Value2 = 3922.0908
$ws.Range("I85").Value2 = 4250.4443
$ws.Range("J85").Value2 = 2444.5
$ws.Range("K85").Value2 = 4250.4443
$ws.Range("L85").Value2 = 2444.5
$ws.Range("M85").Value2 = -3002.4443
$ws.Range("N85").Value2 = -4940.5
$ws.Range("H93").Value2 = 5000.8
$ws.Range("I93").Value2 = 5776
$ws.Range("J93").Value2 = 4225.6
$ws.Range("K93").Value2 = 5776
$ws.Range("L93").Value2 = 4225.6
$ws.Range("M93").Value2 = -4528
$ws.Range("N93").Value2 = -6721.6
$ws.Range("H136").Value2 = 2584.7307
$ws.Range("I136").Value2 = 2012.6471
$ws.Range("K136").Value2 = 6037.9413
$ws.Range("M136").Value2 = -3487.9413

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value2 = 0
$ws.Range("J81").Value2 = 20010100
$ws.Range("K81").Value2 = 0
$ws.Range("L81").Value2 = 40020200
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value2 = -40022322
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = 20010100
$ws.Range("K84").Value2 = 0
$ws.Range("L84").Value2 = 200101000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value2 = -200111608
$ws.Range("H107").Value2 = 494.54166
$ws.Range("I107").Value2 = 508.26315
$ws.Range("K107").Value2 = 1524.78945
$ws.Range("M107").Value2 = 395.21055
$ws.Range("H136").Value2 = 5484.5713
$ws.Range("I136").Value2 = 2758.2
$ws.Range("K136").Value2 = 8274.599999999999
$ws.Range("M136").Value2 = -5724.599999999999
